$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.435.58"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.38%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.865.65"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.18%  "

$ws.Range("E4").Value = "  +0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.72"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("E6").Value = "  +0.01%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4781"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2788"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.57%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06534"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.55%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.861.47"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07446"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "16.20"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.058"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "86.88"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6391"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "30.425.97"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.43%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "12.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.31%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "233.67"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +6.39%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007462"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.42%  "

$ws.Range("E21").Value = "  -0.03%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.119"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.77%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.089"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.09%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "168.74"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.300"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "18.14"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.54%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.894"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.72%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.1055"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +12.57%  "

$ws.Range("E29").Value = "  -5.09%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.267"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.966"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.11%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.04977"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.59%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.168"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.90%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7391"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.714"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01943"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +7.10%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.636"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.9125"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.036"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -2.08%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "106.36"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "

$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("E43").Value = "  -2.12%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.562"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -5.81%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "7.156"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.91%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "61.58"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.58%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1222"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.76%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.881"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "33.52"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "

$ws.Range("E50").Value = "  -4.06%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05637"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "
